# Trade #103 closed at 2026-02-17 09:18:17 - unknown UNKNOWN +0.000%

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Summary": update aggregate metrics to reflect the new trade
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1200.05   # Current Capital
$wsSummary.Range("B4").Value = 0.06      # Total P&L $
$wsSummary.Range("B6").Value = 103       # Total Trades
$wsSummary.Range("B7").Value = 43        # Winning Trades
$wsSummary.Range("B9").Value = 41.75     # Win Rate %

# ---------------------------------------------------------------------
# Sheet "Strategy Status": update the MarketMaking strategy row (row 4)
# ---------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 100.05     # Capital
$wsStatus.Range("D4").Value = 103        # Trades
$wsStatus.Range("E4").Value = 0.06       # P&L $
$wsStatus.Range("F4").Value = 0.05       # P&L %
$wsStatus.Range("G4").Value = 41.75      # Win Rate %

# ---------------------------------------------------------------------
# Sheet "All Trades": append new trade #103 as row 104
# ---------------------------------------------------------------------
$wsAllTrades = $wb.Worksheets.Item("All Trades")
$wsAllTrades.Range("A104").Value = 103
$wsAllTrades.Range("B104").NumberFormat = "@"
$wsAllTrades.Range("B104").Value = "2026-02-17"
$wsAllTrades.Range("C104").Value = "09:18:11"
$wsAllTrades.Range("D104").Value = "MarketMaking"
$wsAllTrades.Range("E104").Value = "DOWN"
$wsAllTrades.Range("F104").Value = 0.76
$wsAllTrades.Range("G104").Value = 0.79
$wsAllTrades.Range("H104").Value = "CLOSED"
$wsAllTrades.Range("I104").Value = 3.9474
$wsAllTrades.Range("J104").Value = 0.03
$wsAllTrades.Range("K104").Value = 100.05
$wsAllTrades.Range("L104").Value = 0
$wsAllTrades.Range("M104").Value = 0
$wsAllTrades.Range("N104").Value = 0.6
$wsAllTrades.Range("O104").Value = "Normal spread capture: 19600 bps"
$wsAllTrades.Range("P104").Value = "early_exit"
$wsAllTrades.Range("Q104").Value = 0.13

# ---------------------------------------------------------------------
# Sheet "MarketMaking": append the same new trade #103 as row 104
# ---------------------------------------------------------------------
$wsMarketMaking = $wb.Worksheets.Item("MarketMaking")
$wsMarketMaking.Range("A104").Value = 103
$wsMarketMaking.Range("B104").NumberFormat = "@"
$wsMarketMaking.Range("B104").Value = "2026-02-17"
$wsMarketMaking.Range("C104").Value = "09:18:11"
$wsMarketMaking.Range("D104").Value = "MarketMaking"
$wsMarketMaking.Range("E104").Value = "DOWN"
$wsMarketMaking.Range("F104").Value = 0.76
$wsMarketMaking.Range("G104").Value = 0.79
$wsMarketMaking.Range("H104").Value = "CLOSED"
$wsMarketMaking.Range("I104").Value = 3.9474
$wsMarketMaking.Range("J104").Value = 0.03
$wsMarketMaking.Range("K104").Value = 100.05
$wsMarketMaking.Range("L104").Value = 0
$wsMarketMaking.Range("M104").Value = 0
$wsMarketMaking.Range("N104").Value = 0.6
$wsMarketMaking.Range("O104").Value = "Normal spread capture: 19600 bps"
$wsMarketMaking.Range("P104").Value = "early_exit"
$wsMarketMaking.Range("Q104").Value = 0.13
